$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2:G13").NumberFormat = "@"
$ws.Range("B2").Value = "0.06"
$ws.Range("C2").Value = "-0.11"
$ws.Range("D2").Value = "0.03"
$ws.Range("E2").Value = "-0.08"
$ws.Range("F2").Value = "0.07"
$ws.Range("G2").Value = "-0.03"
$ws.Range("B3").Value = "0.12"
$ws.Range("C3").Value = "-0.02"
$ws.Range("D3").Value = "0.09"
$ws.Range("E3").Value = "0.07"
$ws.Range("F3").Value = "0.04"
$ws.Range("G3").Value = "-0.2*"
$ws.Range("B4").Value = "0.09"
$ws.Range("C4").Value = "-0.04"
$ws.Range("D4").Value = "0.14"
$ws.Range("E4").Value = "0.19"
$ws.Range("F4").Value = "0.12"
$ws.Range("G4").Value = "-0.19"
$ws.Range("B5").Value = "0.14"
$ws.Range("C5").Value = "-0.1"
$ws.Range("D5").Value = "0.18"
$ws.Range("E5").Value = "0.04"
$ws.Range("F5").Value = "0.03"
$ws.Range("G5").Value = "-0.24*"
$ws.Range("B6").Value = "0.01"
$ws.Range("C6").Value = "-0.19"
$ws.Range("D6").Value = "-0.07"
$ws.Range("E6").Value = "-0.17"
$ws.Range("F6").Value = "0.14"
$ws.Range("G6").Value = "0.02"
$ws.Range("B7").Value = "0.06"
$ws.Range("C7").Value = "-0.29**"
$ws.Range("D7").Value = "0.03"
$ws.Range("E7").Value = "-0.33***"
$ws.Range("F7").Value = "0.05"
$ws.Range("G7").Value = "0.01"
$ws.Range("B8").Value = "-0.0"
$ws.Range("C8").Value = "-0.23*"
$ws.Range("D8").Value = "-0.03"
$ws.Range("E8").Value = "-0.26**"
$ws.Range("F8").Value = "-0.19"
$ws.Range("G8").Value = "-0.02"
$ws.Range("B9").Value = "0.04"
$ws.Range("C9").Value = "-0.17"
$ws.Range("D9").Value = "0.08"
$ws.Range("E9").Value = "-0.24*"
$ws.Range("F9").Value = "-0.25*"
$ws.Range("G9").Value = "-0.08"
$ws.Range("B10").Value = "-0.07"
$ws.Range("C10").Value = "-0.12"
$ws.Range("D10").Value = "-0.03"
$ws.Range("E10").Value = "-0.16"
$ws.Range("F10").Value = "-0.15"
$ws.Range("G10").Value = "-0.07"
$ws.Range("B11").Value = "-0.15"
$ws.Range("C11").Value = "-0.19"
$ws.Range("D11").Value = "-0.08"
$ws.Range("E11").Value = "-0.18"
$ws.Range("F11").Value = "0.02"
$ws.Range("G11").Value = "0.03"
$ws.Range("B12").Value = "-0.23*"
$ws.Range("C12").Value = "-0.17"
$ws.Range("D12").Value = "-0.05"
$ws.Range("E12").Value = "-0.08"
$ws.Range("F12").Value = "0.01"
$ws.Range("G12").Value = "-0.08"
$ws.Range("B13").Value = "-0.11"
$ws.Range("C13").Value = "0.06"
$ws.Range("D13").Value = "0.06"
$ws.Range("E13").Value = "0.04"
$ws.Range("F13").Value = "-0.18"
$ws.Range("G13").Value = "0.02"
$ws.Range("B2:G13").Style = "Normal"
